# Auto-generated edits applying numeric corrections to profit/cost columns (H:N)
# across multiple sheets, per the scheduled runner update.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 7131.8667
$ws.Range("I19").Value = 396.30768
$ws.Range("J19").Value = 12282.588
$ws.Range("K19").Value = 396.30768
$ws.Range("L19").Value = 12282.588
$ws.Range("M19").Value = -221.30768
$ws.Range("N19").Value = -12632.588

# Row 33
$ws.Range("H33").Value = 506
$ws.Range("I33").Value = 575.8570999999999
$ws.Range("K33").Value = 575.8570999999999
$ws.Range("M33").Value = -346.8570999999999

# Row 132
$ws.Range("H132").Value = 4320.0835
$ws.Range("I132").Value = 4969
$ws.Range("J132").Value = 3022.25
$ws.Range("K132").Value = 14907
$ws.Range("L132").Value = 9066.75
$ws.Range("M132").Value = -12377
$ws.Range("N132").Value = -14126.75

# Row 137
$ws.Range("H137").Value = 1661.7354
$ws.Range("I137").Value = 1556.1578
$ws.Range("J137").Value = 1795.4667
$ws.Range("K137").Value = 4668.4734
$ws.Range("L137").Value = 5386.4001
$ws.Range("M137").Value = -2118.4734
$ws.Range("N137").Value = -10486.4001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 110
$ws.Range("H110").Value = 654.05554
$ws.Range("I110").Value = 535.8125
$ws.Range("J110").Value = 1600
$ws.Range("K110").Value = 535.8125
$ws.Range("L110").Value = 1600
$ws.Range("M110").Value = 1509.1875
$ws.Range("N110").Value = -5690

# Row 122
$ws.Range("H122").Value = 61698.35
$ws.Range("I122").Value = 167948
$ws.Range("J122").Value = 3744
$ws.Range("K122").Value = 503844
$ws.Range("L122").Value = 11232
$ws.Range("M122").Value = -501394
$ws.Range("N122").Value = -16132

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 12
$ws.Range("H12").Value = 537.375
$ws.Range("I12").Value = 519.8
$ws.Range("J12").Value = 566.6667
$ws.Range("K12").Value = 519.8
$ws.Range("L12").Value = 566.6667
$ws.Range("M12").Value = -351.8
$ws.Range("N12").Value = -902.6667

# Row 59
$ws.Range("H59").Value = 60780
$ws.Range("J59").Value = 60780
$ws.Range("L59").Value = 60780
$ws.Range("N59").Value = -62474

# Row 86
$ws.Range("H86").Value = 1958.1538
$ws.Range("I86").Value = 1795.8422
$ws.Range("J86").Value = 2398.7144
$ws.Range("K86").Value = 1795.8422
$ws.Range("L86").Value = 2398.7144
$ws.Range("M86").Value = -672.8422
$ws.Range("N86").Value = -4644.7144

# Row 89
$ws.Range("H89").Value = 1958.1538
$ws.Range("I89").Value = 1795.8422
$ws.Range("J89").Value = 2398.7144
$ws.Range("K89").Value = 8979.210999999999
$ws.Range("L89").Value = 11993.572
$ws.Range("M89").Value = -3363.210999999999
$ws.Range("N89").Value = -23225.572

# Row 94
$ws.Range("H94").Value = 500750
$ws.Range("I94").Value = 1000000
$ws.Range("J94").Value = 1500
$ws.Range("K94").Value = 1000000
$ws.Range("L94").Value = 1500
$ws.Range("M94").Value = -999549
$ws.Range("N94").Value = -2402

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5022.7144
$ws.Range("I31").Value = 1357.2
$ws.Range("J31").Value = 6727.6045
$ws.Range("K31").Value = 1357.2
$ws.Range("L31").Value = 6727.6045
$ws.Range("M31").Value = -1062.2
$ws.Range("N31").Value = -7317.6045

# Row 34
$ws.Range("H34").Value = 5022.7144
$ws.Range("I34").Value = 1357.2
$ws.Range("J34").Value = 6727.6045
$ws.Range("K34").Value = 1357.2
$ws.Range("L34").Value = 6727.6045
$ws.Range("M34").Value = -1155.2
$ws.Range("N34").Value = -7131.6045

# Row 58
$ws.Range("H58").Value = 2621.739
$ws.Range("I58").Value = 2433.4375
$ws.Range("J58").Value = 3052.1428
$ws.Range("K58").Value = 2433.4375
$ws.Range("L58").Value = 3052.1428
$ws.Range("M58").Value = -2230.4375
$ws.Range("N58").Value = -3458.1428

# Row 132
$ws.Range("H132").Value = 2672.652
$ws.Range("I132").Value = 2447.0715
$ws.Range("J132").Value = 3023.5557
$ws.Range("K132").Value = 7341.2145
$ws.Range("L132").Value = 9070.667099999999
$ws.Range("M132").Value = -4811.2145
$ws.Range("N132").Value = -14130.6671

# Row 134
$ws.Range("H134").Value = 10006742
$ws.Range("I134").Value = 14713660
$ws.Range("J134").Value = 4541.875
$ws.Range("K134").Value = 44140980
$ws.Range("L134").Value = 13625.625
$ws.Range("M134").Value = -44138445
$ws.Range("N134").Value = -18695.625

# Row 136
$ws.Range("H136").Value = 2621.739
$ws.Range("I136").Value = 2433.4375
$ws.Range("J136").Value = 3052.1428
$ws.Range("K136").Value = 7300.3125
$ws.Range("L136").Value = 9156.428400000001
$ws.Range("M136").Value = -4750.3125
$ws.Range("N136").Value = -14256.4284

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 688.7547
$ws.Range("I5").Value = 522.6129
$ws.Range("J5").Value = 922.86365
$ws.Range("K5").Value = 1567.8387
$ws.Range("L5").Value = 2768.59095
$ws.Range("M5").Value = -1455.8387
$ws.Range("N5").Value = -2992.59095

# Row 49
$ws.Range("H49").Value = 6950.923
$ws.Range("J49").Value = 6950.923
$ws.Range("L49").Value = 20852.769
$ws.Range("N49").Value = -21164.769

# Row 122
$ws.Range("H122").Value = 2323.2034
$ws.Range("I122").Value = 418.12
$ws.Range("J122").Value = 3724
$ws.Range("K122").Value = 3763.08
$ws.Range("L122").Value = 33516
$ws.Range("M122").Value = -1313.08
$ws.Range("N122").Value = -38416

# Row 131
$ws.Range("H131").Value = 4649.3
$ws.Range("J131").Value = 7901.706
$ws.Range("L131").Value = 23705.118
$ws.Range("N131").Value = -33785.118

# Row 135
$ws.Range("H135").Value = 688.7547
$ws.Range("I135").Value = 522.6129
$ws.Range("J135").Value = 922.86365
$ws.Range("K135").Value = 4703.5161
$ws.Range("L135").Value = 8305.772849999999
$ws.Range("M135").Value = -2168.5161
$ws.Range("N135").Value = -13375.77285

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 3340
$ws.Range("I126").Value = 2900
$ws.Range("J126").Value = 3633.3333
$ws.Range("K126").Value = 8700
$ws.Range("L126").Value = 10899.9999
$ws.Range("M126").Value = -6230
$ws.Range("N126").Value = -15839.9999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 2742.25
$ws.Range("I82").Value = 2967.25
$ws.Range("J82").Value = 2292.25
$ws.Range("K82").Value = 2967.25
$ws.Range("L82").Value = 2292.25
$ws.Range("M82").Value = -2606.25
$ws.Range("N82").Value = -3014.25

# Row 85
$ws.Range("H85").Value = 2742.25
$ws.Range("I85").Value = 2967.25
$ws.Range("J85").Value = 2292.25
$ws.Range("K85").Value = 2967.25
$ws.Range("L85").Value = 2292.25
$ws.Range("M85").Value = -1719.25
$ws.Range("N85").Value = -4788.25

# Row 93
$ws.Range("H93").Value = 22020
$ws.Range("I93").Value = 26275
$ws.Range("K93").Value = 26275
$ws.Range("M93").Value = -25027

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 2644
$ws.Range("I113").Value = 3106.6667
$ws.Range("J113").Value = 1950
$ws.Range("K113").Value = 9320.000100000001
$ws.Range("L113").Value = 5850
$ws.Range("M113").Value = -7150.000100000001
$ws.Range("N113").Value = -10190

# Row 132
$ws.Range("H132").Value = 4945619.5
$ws.Range("I132").Value = 1882.9286
$ws.Range("J132").Value = 17159556
$ws.Range("K132").Value = 5648.7858
$ws.Range("L132").Value = 51478668
$ws.Range("M132").Value = -3118.7858
$ws.Range("N132").Value = -51483728

# Row 136
$ws.Range("H136").Value = 4892.2593
$ws.Range("I136").Value = 4572.875
$ws.Range("K136").Value = 13718.625
$ws.Range("M136").Value = -11168.625
